$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

function Clear-CellValue($ws, $ref) {
    $ws.Range($ref).ClearContents()
}


# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws "H5" 122.666664
Set-CellValue $ws "I5" 124
Set-CellValue $ws "K5" 124
Set-CellValue $ws "M5" -9
Set-CellValue $ws "H6" 232.25
Set-CellValue $ws "I6" 232.25
Set-CellValue $ws "J6" 0
Set-CellValue $ws "K6" 696.75
Set-CellValue $ws "L6" 0
Set-CellValue $ws "M6" -584.75
Clear-CellValue $ws "N6"
Set-CellValue $ws "H8" 15
Set-CellValue $ws "I8" 10
Set-CellValue $ws "J8" 20
Set-CellValue $ws "K8" 30
Set-CellValue $ws "L8" 60
Set-CellValue $ws "M8" 109
Set-CellValue $ws "N8" -338
Set-CellValue $ws "H21" 711.3333
Set-CellValue $ws "I21" 711.3333
Set-CellValue $ws "K21" 711.3333
Set-CellValue $ws "M21" -243.3333
Set-CellValue $ws "H23" 711.3333
Set-CellValue $ws "I23" 711.3333
Set-CellValue $ws "K23" 711.3333
Set-CellValue $ws "M23" -477.3333
Set-CellValue $ws "H32" 999.75
Set-CellValue $ws "I32" 999.75
Set-CellValue $ws "K32" 999.75
Set-CellValue $ws "M32" -673.75
Set-CellValue $ws "H112" 1180.4
Set-CellValue $ws "I112" 1050
Set-CellValue $ws "J112" 1236.2858
Set-CellValue $ws "K112" 3150
Set-CellValue $ws "L112" 3708.8574
Set-CellValue $ws "M112" -2042
Set-CellValue $ws "N112" -5924.857400000001
Set-CellValue $ws "H138" 3522.0667
Set-CellValue $ws "J138" 3828.7
Set-CellValue $ws "L138" 11486.1
Set-CellValue $ws "N138" -21766.1

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws "H32" 4478.768
Set-CellValue $ws "I32" 3086.0815
Set-CellValue $ws "K32" 3086.0815
Set-CellValue $ws "M32" -2799.0815
Set-CellValue $ws "H39" 7000
Set-CellValue $ws "I39" 7000
Set-CellValue $ws "K39" 7000
Set-CellValue $ws "M39" -6480
Set-CellValue $ws "H61" 2028.4117
Set-CellValue $ws "I61" 1984.8
Set-CellValue $ws "K61" 1984.8
Set-CellValue $ws "M61" -1772.8
Set-CellValue $ws "H122" 667984.4
Set-CellValue $ws "I122" 715586.9
Set-CellValue $ws "K122" 2146760.7
Set-CellValue $ws "M122" -2144310.7
Set-CellValue $ws "H132" 1456.9697
Set-CellValue $ws "I132" 1534.7931
Set-CellValue $ws "J132" 892.75
Set-CellValue $ws "K132" 4604.379300000001
Set-CellValue $ws "L132" 2678.25
Set-CellValue $ws "M132" -2074.379300000001
Set-CellValue $ws "N132" -7738.25
Set-CellValue $ws "H136" 2028.4117
Set-CellValue $ws "I136" 1984.8
Set-CellValue $ws "K136" 5954.4
Set-CellValue $ws "M136" -3404.4

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws "H24" 4933
Set-CellValue $ws "I24" 4933
Set-CellValue $ws "K24" 4933
Set-CellValue $ws "M24" -4698
Set-CellValue $ws "H80" 458.16666
Set-CellValue $ws "I80" 433.33334
Set-CellValue $ws "K80" 433.33334
Set-CellValue $ws "M80" 564.66666
Set-CellValue $ws "H83" 458.16666
Set-CellValue $ws "I83" 433.33334
Set-CellValue $ws "K83" 2166.6667
Set-CellValue $ws "M83" 2825.3333
Set-CellValue $ws "H86" 4926.5
Set-CellValue $ws "I86" 4999
Set-CellValue $ws "J86" 4890.25
Set-CellValue $ws "K86" 4999
Set-CellValue $ws "L86" 4890.25
Set-CellValue $ws "M86" -3876
Set-CellValue $ws "N86" -7136.25
Set-CellValue $ws "H89" 4926.5
Set-CellValue $ws "I89" 4999
Set-CellValue $ws "J89" 4890.25
Set-CellValue $ws "K89" 24995
Set-CellValue $ws "L89" 24451.25
Set-CellValue $ws "M89" -19379
Set-CellValue $ws "N89" -35683.25
Set-CellValue $ws "H94" 672.5
Set-CellValue $ws "J94" 129
Set-CellValue $ws "L94" 129
Set-CellValue $ws "N94" -1031
Set-CellValue $ws "H105" 3053.9443
Set-CellValue $ws "I105" 2557.25
Set-CellValue $ws "J105" 3451.3
Set-CellValue $ws "K105" 2557.25
Set-CellValue $ws "L105" 3451.3
Set-CellValue $ws "M105" -810.25
Set-CellValue $ws "N105" -6945.3
Set-CellValue $ws "H107" 1983.75
Set-CellValue $ws "J107" 1565.6666
Set-CellValue $ws "L107" 1565.6666
Set-CellValue $ws "N107" -5405.6666
Set-CellValue $ws "H134" 2782.5715
Set-CellValue $ws "I134" 1895.6
Set-CellValue $ws "K134" 5686.799999999999
Set-CellValue $ws "M134" -3151.799999999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws "H4" 4220.25
Set-CellValue $ws "I4" 2900
Set-CellValue $ws "J4" 4660.3335
Set-CellValue $ws "K4" 2900
Set-CellValue $ws "L4" 4660.3335
Set-CellValue $ws "N4" -4884.3335
Set-CellValue $ws "M4" -2788
Set-CellValue $ws "H12" 5201.2856
Set-CellValue $ws "I12" 3901.5
Set-CellValue $ws "J12" 13000
Set-CellValue $ws "K12" 3901.5
Set-CellValue $ws "L12" 13000
Set-CellValue $ws "M12" -3731.5
Set-CellValue $ws "N12" -13340
Set-CellValue $ws "H105" 3189.5417
Set-CellValue $ws "I105" 855.1
Set-CellValue $ws "J105" 4857
Set-CellValue $ws "K105" 855.1
Set-CellValue $ws "L105" 4857
Set-CellValue $ws "M105" 891.9
Set-CellValue $ws "N105" -8351
Set-CellValue $ws "H122" 1999
Set-CellValue $ws "I122" 1999
Set-CellValue $ws "K122" 5997
Set-CellValue $ws "M122" -3547

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws "H4" 2341456.8
Set-CellValue $ws "I4" 2341456.8
Set-CellValue $ws "J4" 0
Set-CellValue $ws "K4" 7024370.399999999
Set-CellValue $ws "L4" 0
Set-CellValue $ws "M4" -7024258.399999999
Clear-CellValue $ws "N4"
Set-CellValue $ws "H113" 799.5
Set-CellValue $ws "I113" 799.5
Set-CellValue $ws "K113" 2398.5
Set-CellValue $ws "M113" -228.5

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws "H20" 44444
Set-CellValue $ws "J20" 44444
Set-CellValue $ws "L20" 44444
Set-CellValue $ws "N20" -44934
Set-CellValue $ws "H24" 31271.572
Set-CellValue $ws "J24" 31271.572
Set-CellValue $ws "L24" 31271.572
Set-CellValue $ws "N24" -31617.572
Set-CellValue $ws "H70" 6833.3335
Set-CellValue $ws "I70" 0
Set-CellValue $ws "J70" 6833.3335
Set-CellValue $ws "K70" 0
Set-CellValue $ws "L70" 6833.3335
Clear-CellValue $ws "M70"
Set-CellValue $ws "N70" -7373.3335
Set-CellValue $ws "H73" 6833.3335
Set-CellValue $ws "I73" 0
Set-CellValue $ws "J73" 6833.3335
Set-CellValue $ws "K73" 0
Set-CellValue $ws "L73" 6833.3335
Clear-CellValue $ws "M73"
Set-CellValue $ws "N73" -8705.333500000001
Set-CellValue $ws "H92" 13993.6
Set-CellValue $ws "I92" 10000
Set-CellValue $ws "J92" 14992
Set-CellValue $ws "K92" 10000
Set-CellValue $ws "L92" 14992
Set-CellValue $ws "M92" -8128
Set-CellValue $ws "N92" -18736
Set-CellValue $ws "H132" 2543.348
Set-CellValue $ws "I132" 1909.1875
Set-CellValue $ws "K132" 5727.5625
Set-CellValue $ws "M132" -3197.5625
Set-CellValue $ws "H134" 84853.57000000001
Set-CellValue $ws "J134" 84853.57000000001
Set-CellValue $ws "L134" 254560.71
Set-CellValue $ws "N134" -259630.71

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws "H53" 13999.667
Set-CellValue $ws "I53" 13999.667
Set-CellValue $ws "K53" 13999.667
Set-CellValue $ws "M53" -13481.667
Set-CellValue $ws "H93" 1133.1765
Set-CellValue $ws "I93" 920.61536
Set-CellValue $ws "J93" 1824
Set-CellValue $ws "K93" 920.61536
Set-CellValue $ws "L93" 1824
Set-CellValue $ws "M93" 327.38464
Set-CellValue $ws "N93" -4320

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws "H18" 13502.75
Set-CellValue $ws "J18" 13502.75
Set-CellValue $ws "L18" 13502.75
Set-CellValue $ws "N18" -13848.75
Set-CellValue $ws "H29" 18000
Set-CellValue $ws "J29" 18000
Set-CellValue $ws "L29" 18000
Set-CellValue $ws "N29" -18580
Set-CellValue $ws "H62" 6699.6665
Set-CellValue $ws "J62" 7014.0713
Set-CellValue $ws "L62" 7014.0713
Set-CellValue $ws "N62" -8262.0713
Set-CellValue $ws "H65" 6699.6665
Set-CellValue $ws "J65" 7014.0713
Set-CellValue $ws "L65" 35070.35649999999
Set-CellValue $ws "N65" -41310.35649999999
Set-CellValue $ws "H132" 3224.4707
Set-CellValue $ws "I132" 2562.7693
Set-CellValue $ws "K132" 7688.3079
Set-CellValue $ws "M132" -5158.3079

Write-Output "Done applying Seraph_Profits updates"
